# Add the new "synonyms" / "GF_synonyms" column (H) to Sheet1, and move the
# current selection the way it ended up after the edit (cell I3), matching
# the author's recorded change ("Committing the changes done in the tests
# to run them on IE 11").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header cell (H1) + new data cell (H2) for the extra "synonyms" column.
$ws.Cells.Item(1, 8).Value = "synonyms"
$ws.Cells.Item(2, 8).Value = "GF_synonyms"

# Leave the selection where the author's session left it.
$ws.Range("I3").Select() | Out-Null
